# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sorts / refreshes the "Estado de Cuenta" detail table (rows 16-29) so
# that the two workers' mora periods run in ascending order (2405..2411),
# interleaving EVER EDUARDO CARDENAS DE LA OSSA (CC 1053004704) and
# FREDY HUERTAS LOPEZ (CC 74339257) for each period, instead of the old
# per-worker / descending-period grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, TipoDoc(B), NumDoc(C), Nombre(D), PeriodoMora(E), ValorMora(F), SalarioBasico(G)
$rows = @(
    @(16, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2405", 58720,  1468000),
    @(17, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2405", 94286,  2357150),
    @(18, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2406", 58720,  1468000),
    @(19, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2406", 94286,  2357150),
    @(20, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2407", 58720,  1468000),
    @(21, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2407", 94286,  2357150),
    @(22, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2408", 58720,  1468000),
    @(23, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2408", 94286,  2357150),
    @(24, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2409", 58720,  1468000),
    @(25, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2409", 94286,  2357150),
    @(26, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2410", 58720,  1468000),
    @(27, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2410", 94286,  2357150),
    @(28, "CC", "1053004704", "EVER EDUARDO CARDENAS DE LA OSSA", "2411", 27402,  1468000),
    @(29, "CC", "74339257",   "FREDY HUERTAS LOPEZ",              "2411", 44000,  2357150)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($rowNum, 3).Value = $r[2]   # C: N Doc Trabajador
    $ws.Cells.Item($rowNum, 4).Value = $r[3]   # D: Nombre Trabajador
    $ws.Cells.Item($rowNum, 5).Value = $r[4]   # E: Periodo Mora
    $ws.Cells.Item($rowNum, 6).Value = $r[5]   # F: Valor Mora
    $ws.Cells.Item($rowNum, 7).Value = $r[6]   # G: Salario Basico
}
